$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two cells in row 3 (2025-01-22) for greenpeace.eu (H) and wwfeu.bsky.social (I)
# were removed/cleared - their bot accounts posted 0 (no data) that day.
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()

# Restore the author's last selection/view state (cell J3 selected)
[void]$ws.Range("J3").Select()
